# Add a new row of daily-expense data (row 27) below the existing data,
# matching the date-style formatting used by the rows above it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date cell's format (short-date numFmt) from the prior row, then
# overwrite the value - avoids minting a brand-new custom number format.
$ws.Range("A26").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A27").Value = 43816

$ws.Range("B27").Value = 2
$ws.Range("C27").Value = 0
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 12.5
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1
$ws.Range("L27").Value = 1
$ws.Range("M27").Value = 3

# Move the active selection to K27, matching the post-edit cursor position.
$ws.Range("K27").Select()
